$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Reference to the existing "Message" sheet - used as the template for the
# two new "*_Message" sheets (they duplicate its 10-row connection blurb).
# ---------------------------------------------------------------------------
$msgSheet = $wb.Worksheets.Item("Message")

# ---------------------------------------------------------------------------
# 1) CypherOutput_Message  -  verbatim copy of "Message", inserted right
#    after it.
# ---------------------------------------------------------------------------
$msgSheet.Copy($null, $msgSheet)
$cypherMsgSheet = $wb.Worksheets.Item($msgSheet.Index + 1)
$cypherMsgSheet.Name = "CypherOutput_Message"

# ---------------------------------------------------------------------------
# 2) StatOutput  -  new sheet with the aggregate counts, inserted after
#    CypherOutput_Message.
# ---------------------------------------------------------------------------
$statSheet = $wb.Worksheets.Add($null, $cypherMsgSheet)
$statSheet.Name = "StatOutput"

$headerRange = $statSheet.Range("A1:D1")
$headerRange.Cells.Item(1, 1).Formula = "=""number_of_files"""
$headerRange.Cells.Item(1, 2).Formula = "=""number_of_sample"""
$headerRange.Cells.Item(1, 3).Formula = "=""number_of_cases"""
$headerRange.Cells.Item(1, 4).Formula = "=""number_of_study"""
$headerRange.Copy()
$headerRange.PasteSpecial(-4163)

$dataRange = $statSheet.Range("A2:D2")
$dataRange.Cells.Item(1, 1).Formula = "=""0"""
$dataRange.Cells.Item(1, 2).Formula = "=""0"""
$dataRange.Cells.Item(1, 3).Formula = "=""7"""
$dataRange.Cells.Item(1, 4).Formula = "=""1"""
$dataRange.Copy()
$dataRange.PasteSpecial(-4163)

# ---------------------------------------------------------------------------
# 3) StatOutput_Message  -  copy of "Message" (rows 1-10) with the same
#    10-row block repeated again (rows 11-20), except the Cypher-text row
#    (row 18) is replaced with the new StatOutput query text.
# ---------------------------------------------------------------------------
$msgSheet.Copy($null, $statSheet)
$statMsgSheet = $wb.Worksheets.Item($statSheet.Index + 1)
$statMsgSheet.Name = "StatOutput_Message"

$firstBlock = $statMsgSheet.Range("A1:A10")
$firstBlock.Copy()
$secondBlock = $statMsgSheet.Range("A11:A20")
$secondBlock.PasteSpecial(-4163)

$statQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN['Stage 4']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
$statMsgSheet.Range("A18").Value = $statQuery

# Restore selection/active sheet back to the first sheet, matching the
# original workbook's "tabSelected" sheet.
$wb.Worksheets.Item("CypherOutput").Activate()
